# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "periodos" (period) table on the sheet grows from 3 rows to 4 rows:
# old periods 2504/2505/2506 are replaced by 2507/2506/2505/2504 (a new,
# more recent period 2507 is added, and an older period 2504 is added back
# at the bottom), and the totals (Valor Mora / Cant. Periodos) are updated
# to reflect the extra period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row just above the last existing period row (row 18).
# This pushes the former row 18 down to row 19, and the signature block
# further down (rows 23/24 -> 24/25) automatically, along with the merged
# cell ranges.
$ws.Rows(18).Insert()

# The newly inserted row 18 has no formatting yet - copy the look of the
# row above it (row 17, one of the "middle" period rows) into it.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new row 18 with the same worker info as the other period
# rows, for period 2505.
$ws.Range("B18").Value2 = $ws.Range("B17").Value2
$ws.Range("C18").Value2 = $ws.Range("C17").Value2
$ws.Range("D18").Value2 = $ws.Range("D17").Value2
$ws.Range("E18").Value2 = "2505"
$ws.Range("F18").Value2 = 46400
$ws.Range("G18").Value2 = 1160000

# Shift the remaining periods: newest period (2507) goes on top, followed
# by 2506, 2505 (new row 18), and 2504 at the bottom (row 19, which kept
# the old row 18's formatting/values except for the period code).
$ws.Range("E16").Value2 = "2507"
$ws.Range("E17").Value2 = "2506"
$ws.Range("E19").Value2 = "2504"

# Update the summary totals for the extra period.
$ws.Range("E11").Value2 = 185600
$ws.Range("F13").Value2 = 4
